$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 354.89795
$ws.Range("J17").Value = 354.89795
$ws.Range("L17").Value = 1064.69385
$ws.Range("N17").Value = -1400.69385
$ws.Range("H40").Value = 10487.167
$ws.Range("J40").Value = 12364.6
$ws.Range("L40").Value = 12364.6
$ws.Range("N40").Value = -12714.6
$ws.Range("H64").Value = 5135.4116
$ws.Range("I64").Value = 4609.9
$ws.Range("J64").Value = 5886.143
$ws.Range("K64").Value = 4609.9
$ws.Range("L64").Value = 5886.143
$ws.Range("M64").Value = -4361.9
$ws.Range("N64").Value = -6382.143
$ws.Range("H67").Value = 5135.4116
$ws.Range("I67").Value = 4609.9
$ws.Range("J67").Value = 5886.143
$ws.Range("K67").Value = 4609.9
$ws.Range("L67").Value = 5886.143
$ws.Range("M67").Value = -3751.9
$ws.Range("N67").Value = -7602.143
$ws.Range("H70").Value = 2854.5483
$ws.Range("I70").Value = 3564.6
$ws.Range("J70").Value = 1563.5454
$ws.Range("K70").Value = 10693.8
$ws.Range("L70").Value = 4690.6362
$ws.Range("M70").Value = -10423.8
$ws.Range("N70").Value = -5230.6362
$ws.Range("H73").Value = 2854.5483
$ws.Range("I73").Value = 3564.6
$ws.Range("J73").Value = 1563.5454
$ws.Range("K73").Value = 10693.8
$ws.Range("L73").Value = 4690.6362
$ws.Range("M73").Value = -9757.799999999999
$ws.Range("N73").Value = -6562.6362
$ws.Range("H98").Value = 628.2632
$ws.Range("I98").Value = 606.06665
$ws.Range("K98").Value = 606.06665
$ws.Range("M98").Value = 891.93335
$ws.Range("H122").Value = 628.2632
$ws.Range("I122").Value = 606.06665
$ws.Range("K122").Value = 1818.19995
$ws.Range("M122").Value = 631.8000500000001
$ws.Range("H129").Value = 796
$ws.Range("I129").Value = 276.5
$ws.Range("J129").Value = 1159.65
$ws.Range("K129").Value = 829.5
$ws.Range("L129").Value = 3478.95
$ws.Range("M129").Value = 4170.5
$ws.Range("N129").Value = -13478.95
$ws.Range("H140").Value = 30545.715
$ws.Range("I140").Value = 24000
$ws.Range("J140").Value = 31636.666
$ws.Range("K140").Value = 24000
$ws.Range("L140").Value = 31636.666
$ws.Range("M140").Value = -18820
$ws.Range("N140").Value = -41996.666
$ws.Range("H141").Value = 2126
$ws.Range("I141").Value = 1764.4445
$ws.Range("J141").Value = 3753
$ws.Range("K141").Value = 5293.333500000001
$ws.Range("L141").Value = 11259
$ws.Range("M141").Value = -113.3335000000006
$ws.Range("N141").Value = -21619

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1466350
$ws.Range("I32").Value = 1705050
$ws.Range("J32").Value = 4311.75
$ws.Range("K32").Value = 1705050
$ws.Range("L32").Value = 4311.75
$ws.Range("M32").Value = -1704763
$ws.Range("N32").Value = -4885.75
$ws.Range("H61").Value = 410875.2
$ws.Range("I61").Value = 314692.8
$ws.Range("J61").Value = 591924.4
$ws.Range("K61").Value = 314692.8
$ws.Range("L61").Value = 591924.4
$ws.Range("M61").Value = -314480.8
$ws.Range("N61").Value = -592348.4
$ws.Range("H136").Value = 410875.2
$ws.Range("I136").Value = 314692.8
$ws.Range("J136").Value = 591924.4
$ws.Range("K136").Value = 944078.3999999999
$ws.Range("L136").Value = 1775773.2
$ws.Range("M136").Value = -941528.3999999999
$ws.Range("N136").Value = -1780873.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3227768.8
$ws.Range("I105").Value = 1887.4
$ws.Range("J105").Value = 16668942
$ws.Range("K105").Value = 1887.4
$ws.Range("L105").Value = 16668942
$ws.Range("M105").Value = -140.4000000000001
$ws.Range("N105").Value = -16672436

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 738.5172
$ws.Range("I16").Value = 756.13336
$ws.Range("J16").Value = 719.6429000000001
$ws.Range("K16").Value = 756.13336
$ws.Range("L16").Value = 719.6429000000001
$ws.Range("M16").Value = -469.13336
$ws.Range("N16").Value = -1293.6429
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("H105").Value = 933.46155
$ws.Range("I105").Value = 880.3488
$ws.Range("J105").Value = 1187.2222
$ws.Range("K105").Value = 880.3488
$ws.Range("L105").Value = 1187.2222
$ws.Range("M105").Value = 866.6512
$ws.Range("N105").Value = -4681.2222
$ws.Range("H113").Value = 738.5172
$ws.Range("I113").Value = 756.13336
$ws.Range("J113").Value = 719.6429000000001
$ws.Range("K113").Value = 756.13336
$ws.Range("L113").Value = 719.6429000000001
$ws.Range("M113").Value = 1413.86664
$ws.Range("N113").Value = -5059.6429
$ws.Range("H134").Value = 1985.9375
$ws.Range("I134").Value = 1141.909
$ws.Range("J134").Value = 3842.8
$ws.Range("K134").Value = 3425.727
$ws.Range("L134").Value = 11528.4
$ws.Range("M134").Value = -890.7270000000003
$ws.Range("N134").Value = -16598.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1507
$ws.Range("I107").Value = 360
$ws.Range("J107").Value = 1889.3334
$ws.Range("K107").Value = 1080
$ws.Range("L107").Value = 5668.0002
$ws.Range("M107").Value = 840
$ws.Range("N107").Value = -9508.0002
$ws.Range("H109").Value = 4993.3
$ws.Range("I109").Value = 5419
$ws.Range("K109").Value = 16257
$ws.Range("M109").Value = -15217
$ws.Range("H113").Value = 10870085
$ws.Range("I113").Value = 16667201
$ws.Range("J113").Value = 494.375
$ws.Range("K113").Value = 50001603
$ws.Range("L113").Value = 1483.125
$ws.Range("M113").Value = -49999433
$ws.Range("N113").Value = -5823.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3380.9858
$ws.Range("I80").Value = 3573.276
$ws.Range("J80").Value = 2523.077
$ws.Range("K80").Value = 3573.276
$ws.Range("L80").Value = 2523.077
$ws.Range("M80").Value = -2575.276
$ws.Range("N80").Value = -4519.077
$ws.Range("H83").Value = 3380.9858
$ws.Range("I83").Value = 3573.276
$ws.Range("J83").Value = 2523.077
$ws.Range("K83").Value = 17866.38
$ws.Range("L83").Value = 12615.385
$ws.Range("M83").Value = -12874.38
$ws.Range("N83").Value = -22599.385
$ws.Range("H113").Value = 1002.2
$ws.Range("I113").Value = 1203.6666
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 1203.6666
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 966.3334
$ws.Range("N113").Value = -5040
$ws.Range("H127").Value = 23490.842
$ws.Range("J127").Value = 23490.842
$ws.Range("L127").Value = 23490.842
$ws.Range("N127").Value = -33410.842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1133.963
$ws.Range("I16").Value = 1180.68
$ws.Range("J16").Value = 550
$ws.Range("K16").Value = 1180.68
$ws.Range("L16").Value = 550
$ws.Range("M16").Value = -1010.68
$ws.Range("N16").Value = -890
